$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.095.62'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.780.27'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.38%  '
$ws.Range("D5").Formula = "'225.33"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").Formula = "'31.75"
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("D11").Formula = "'0.0946"
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '2.037.28'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.782.17'
$ws.Range("E13").Value = '  -0.39%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Formula = "'10.90"
$ws.Range("E14").Value = '  -3.50%  '
$ws.Range("D15").Value = '34.089.43'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").Formula = "'244.74"
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").Value = '0.0₃0785'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").Formula = "'10.92"
$ws.Range("E21").Value = '  +1.99%  '
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("E24").Value = '  -1.12%  '
$ws.Range("D25").Formula = "'161.33"
$ws.Range("E25").Value = '  -0.28%  '
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").Formula = "'0.0515"
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Formula = "'3.69"
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("E34").Value = '  -2.76%  '
$ws.Range("D35").Value = '1.450.47'
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("E36").Value = '  +3.91%  '
$ws.Range("D37").Formula = "'0.651"
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("D41").Formula = "'80.29"
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Formula = "'13.65"
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").Formula = "'0.0517"
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '1.938.84'
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Formula = "'104.03"
$ws.Range("E50").Value = '  -3.35%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0130'
$ws.Range("E51").Value = '  -6.73%  '
